$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.096.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = "'1.891.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'306.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = "'0.9992"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = "'0.5153"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.82%  '
$ws.Range("D8").Value = "'0.3756"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.43%  '
$ws.Range("D9").Value = "'0.07208"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.49%  '
$ws.Range("D10").Value = "'21.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.50%  '
$ws.Range("D11").Value = "'0.9052"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.10%  '
$ws.Range("D12").Value = "'0.07652"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.37%  '
$ws.Range("D13").Value = "'1.887.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = "'94.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.63%  '
$ws.Range("D15").Value = "'5.272"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("D16").Value = "'0.9991"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = "'0.000008477"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = "'14.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.20%  '
$ws.Range("D19").Value = "'0.9991"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = "'27.114.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("D21").Value = "'5.066"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("D22").Value = "'2.117.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.68%  '
$ws.Range("D23").Value = "'10.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.04%  '
$ws.Range("D24").Value = "'6.415"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").Value = "'145.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("D26").Value = "'2.265"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.22%  '
$ws.Range("D27").Value = "'1.769"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("E28").Value = '  +1.14%  '
$ws.Range("D29").Value = "'114.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.31%  '
$ws.Range("D30").Value = "'4.945"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.75%  '
$ws.Range("D31").Value = "'4.838"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.39%  '
$ws.Range("D32").Value = "'0.09174"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.93%  '
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("E34").Value = '  +7.61%  '
$ws.Range("D35").Value = "'0.7841"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.58%  '
$ws.Range("D36").Value = "'2.986"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("D37").Value = "'3.287"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'2.617"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.66%  '
$ws.Range("D39").Value = "'0.02000"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("D40").Value = "'0.5596"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.25%  '
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("D42").Value = "'9.074"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.94%  '
$ws.Range("D43").Value = "'6.628"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.39%  '
$ws.Range("D44").Value = "'117.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("D45").Value = "'0.1510"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.76%  '
$ws.Range("D46").Value = "'0.4800"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.44%  '
$ws.Range("D47").Value = "'10.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("D48").Value = "'0.9990"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("E49").Value = '  +2.06%  '
$ws.Range("D50").Value = "'37.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("D51").Value = "'63.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.58%  '
